$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.179.22"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "3.520.52"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.79%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +4.35%  "
$ws.Range("E9").Value = "  +6.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "4.128.38"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.43%  "
$ws.Range("E15").Value = "  +1.53%  "
$ws.Range("D16").Value = "67.167.77"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "3.480.51"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("E19").Value = "  +1.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "394.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.541"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.52%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("E25").Value = "  -3.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.62%  "
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E29").Value = "  -2.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "162.94"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.898"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.83%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.60%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.26%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "27.54"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.46%  "
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.18%  "
$ws.Range("D44").Value = "2.807.49"
$ws.Range("E44").Value = "  -1.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("E46").Value = "  -2.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "337.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.09%  "
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.62"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.849"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.58%  "
